# Update loading_percent results for the "case with 380 kV done" commit.
# Columns C,D,E,F,G,H,L on rows 2-25 of Sheet1 get new computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in order) whose values change for every data row.
$cols = @("C", "D", "E", "F", "G", "H", "L")

# Each entry: row number followed by the new values for C, D, E, F, G, H, L (in that order).
$data = @(
    @(2, 13.5590018579668, 8.66842573984802, 13.94350578747537, 47.90768081590532, 72.18944697254683, 23.23254646083079, 9.315572113616827),
    @(3, 13.48162486905052, 8.600279589338236, 13.90519190635674, 46.83405872097015, 70.09032292907315, 22.94706004378689, 9.332952975339452),
    @(4, 13.43837151218974, 8.561283301186164, 13.88584101810896, 46.18252198220114, 68.80107669441701, 22.77866971492658, 9.345659042175319),
    @(5, 13.42182111909213, 8.546115299029161, 13.87900244070492, 45.919327259825, 68.27646305094007, 22.71184656557455, 9.351345800526987),
    @(6, 13.41913803068599, 8.543640545595236, 13.87793004497133, 45.87577481330169, 68.18942332164059, 22.70086091555671, 9.352320744637487),
    @(7, 13.43814394651767, 8.561075802476163, 13.88574455525294, 46.17896257887733, 68.7939973141853, 22.77776115841676, 9.345733678552911),
    @(8, 13.5314426164512, 8.644343395871873, 13.9294276621912, 47.53613242312501, 71.46625081769692, 23.13271749797993, 9.321141313725871),
    @(9, 13.74785763166729, 8.829803992955956, 14.04827398022687, 50.24183648997303, 76.66788407378242, 23.88059320178689, 9.289170125078387),
    @(10, 13.92672314136436, 8.978948137811127, 14.15584331716707, 52.23522371739963, 80.42167366441141, 24.45749077919983, 9.275745737809846),
    @(11, 14.01226241521004, 9.049437807733074, 14.20915972431975, 53.13899742026928, 82.10636623692774, 24.72496789563619, 9.271856640215638),
    @(12, 14.04523948718304, 9.076495748221754, 14.22997613552084, 53.48046691368543, 82.7404083124986, 24.82690076231398, 9.270705393871175),
    @(13, 14.03811150536218, 9.07065234905172, 14.22546514605385, 53.40696469294073, 82.60403931071522, 24.80492027918245, 9.270938997608484),
    @(14, 14.01496381920992, 9.051656653789767, 14.21085976804011, 53.16710759304129, 82.1586110395157, 24.73334150394952, 9.27175546923584),
    @(15, 14.00086094517326, 9.040068320617999, 14.20199503714001, 53.02007835300724, 81.88524626135934, 24.68957907039977, 9.272297521850438),
    @(16, 13.92121567229268, 8.974393390303426, 14.15244687154581, 52.17606996575419, 80.31105956290951, 24.44010515494604, 9.276044746176842),
    @(17, 13.87341440701768, 8.934770471476696, 14.12317094580968, 51.65727789742854, 79.33901692144765, 24.2882957251424, 9.278913400493918),
    @(18, 13.84631408193695, 8.912230166716553, 14.10674499607133, 51.35861612822432, 78.77779905419067, 24.20145849695755, 9.280771967287345),
    @(19, 13.8372064263959, 8.904641748294518, 14.10125444054169, 51.25745918716842, 78.58743521373246, 24.17214170793601, 9.281436993052388),
    @(20, 13.87846229037153, 8.938962657424446, 14.12624471867986, 51.71253424990164, 79.4427170769869, 24.30440704848565, 9.278586420053388),
    @(21, 14.02174709709882, 9.057226366662128, 14.21513274800001, 53.23758289492716, 82.28955498714295, 24.75434905562847, 9.271506906930957),
    @(22, 14.11879537796623, 9.136638772420222, 14.27687643601528, 54.22966244685193, 84.12706683342132, 25.05213535015071, 9.26875496511806),
    @(23, 14.06669265417337, 9.094066027348594, 14.2435901042033, 53.70069930277461, 83.14865113589346, 24.89288692482673, 9.270051338267365),
    @(24, 13.876178953662, 8.937066623223808, 14.12485380369222, 51.68755407727876, 79.39584159484076, 24.29712174250689, 9.27873359606016),
    @(25, 13.68577831548531, 8.777318547159329, 14.01256355229882, 49.50739887925343, 75.26962022524258, 23.67314524421309, 9.296062750898884)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
